# Daily attendance processing - reorder the "Recorded By" (column G) author
# lists so "System" sorts first, matching the canonical recorder-name order.
#
# Two exact text substitutions occur, applied to every row of the used range
# whose column-G text currently matches (order-of-authors only changes):
#   "dnasr281@gmail.com, System"               -> "System, dnasr281@gmail.com"
#   "System, backup@backdoor.com, system"      -> "System, system, backup@backdoor.com"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$firstRow = $used.Row
$lastRow = $firstRow + $used.Rows.Count - 1

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = 7
    $val = $cell.Text

    if ($val -eq "dnasr281@gmail.com, System") {
        $cell.Value = "System, dnasr281@gmail.com"
    }
    elseif ($val -eq "System, backup@backdoor.com, system") {
        $cell.Value = "System, system, backup@backdoor.com"
    }
}
